# Add insanity card "Silver Twilight Devotee" and "Aggressive" for Sanctum of Twilight

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the Name/Flavor text first (matches the shared-string ordering of the
# original authoring session), then the B (PlayerCount) numbers, then the
# Rules text, and finally the Clarifications text.

# Names
$ws.Range("A21").Value = "Silver Twilight Devotee"
$ws.Range("C21").Value = "The Order of the Silver Twilight has need of the information and artifacts you uncover."
$ws.Range("A22").Value = "Aggressive"
$ws.Range("C22").Value = "You cannot agree with these people. They insist on wasting precious time and resources. Someone has to advocate for the devil or important matters may be overlooked."

# Rules
$ws.Range("D21").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have 1 or more <i>Bladed Weapons</i> and 1 or more <i>Unique Items.</i> Otherwise, you lose the game.</p>"
$d22 = "<p>You do not win the game as normal. Instead you win only if the investigation is complete and this card has been revealed.</p>`n<p><b>Unchained Rage:</b>At the start of the Investigator Phase, if you are in a space with only one investigator and no other creatures or characters, you may reveal this card and choose that investigator. That investigator suffers one Damage. This may only happen once per game.</p>"
$ws.Range("D22").Value = $d22

# Clarifications
$ws.Range("E22").Value = "It must be just  the two of you together in that space - no other investigators, monsters, other non-player characters, anybody else."

# PlayerCount
$ws.Range("B21").Value = 3
$ws.Range("B22").Value = 2

# Match row heights from the source workbook
$ws.Rows.Item(21).RowHeight = 45
$ws.Rows.Item(22).RowHeight = 90

# Update the selection to match the post-edit view state (the new last row)
$null = $ws.Range("A22").Select()
